# The post originally at row 785 (「遊ぶパンダの光景より可愛いものがありますか？」)
# was removed. Deleting the entire row shifts every following row up by
# one, which matches the diff (old row 786 -> new row 785, ... old row
# 822 -> new row 821) and also updates the sheet's used-range dimension
# from A1:C822 to A1:C821 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(785).Delete()
